# Atualizei dados da bibi
# Update column L (day 11) and column AG (total) for rows 2-6 on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bibi Cell Mundi
$ws.Range("L2").Value = 16650.55
$ws.Range("AG2").Value = 125516.31

# Row 3 - Bibi Cell Vieiralves
$ws.Range("L3").Value = 9429.450000000001
$ws.Range("AG3").Value = 50455.4

# Row 4 - Bibi Cell Ponta Negra
$ws.Range("L4").Value = 3514.36
$ws.Range("AG4").Value = 34784.9

# Row 5 - Bibi Cell Manauara
$ws.Range("L5").Value = 3525
$ws.Range("AG5").Value = 32939.2

# Row 6 - total
$ws.Range("L6").Value = 33119.36
$ws.Range("AG6").Value = 243695.81
